$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "767×5=3835" "738×5=3690"
Replace-Text "959×5=4795" "518×6=3108"
Replace-Text "601×9=5409" "461×3=1383"
Replace-Text "336×6=2016" "501×2=1002"
Replace-Text "596×3=1788" "179×4=716"
Replace-Text "979×8=7832" "405×4=1620"
Replace-Text "583×2=1166" "615×3=1845"
Replace-Text "130×4=520" "266×7=1862"
Replace-Text "527×7=3689" "922×7=6454"
Replace-Text "781×7=5467" "522×6=3132"
Replace-Text "646×2=1292" "527×6=3162"
Replace-Text "455×4=1820" "171×5=855"
Replace-Text "466×9=4194" "269×3=807"
Replace-Text "978×2=1956" "140×4=560"
Replace-Text "984×7=6888" "340×2=680"
Replace-Text "279×9=2511" "965×8=7720"
Replace-Text "466×7=3262" "233×2=466"
Replace-Text "867×8=6936" "950×9=8550"
Replace-Text "220×2=440" "904×3=2712"
Replace-Text "316×7=2212" "445×5=2225"
Replace-Text "613×7=4291" "285×2=570"
Replace-Text "319×3=957" "976×9=8784"
Replace-Text "811×3=2433" "184×8=1472"
Replace-Text "617×9=5553" "926×9=8334"
Replace-Text "267×2=534" "652×9=5868"
